# Apply the "create staging table for DF3" edit:
#  - packages sheet: rename the rd3_portal_release package's label/description
#  - entities sheet: rename the "novelwgs" entity to "rd3_portal_release_freeze3" (Freeze 3)
#  - attributes sheet: repoint every attribute row's "entity" column at the renamed entity

$wb = $excel.ActiveWorkbook

# --- packages sheet -------------------------------------------------------
$packages = $wb.Worksheets.Item("packages")
# Row 3 holds the rd3_portal_release package (A3=rd3_portal_release)
$packages.Cells.Item(3, 2).Value = "Releases"
$packages.Cells.Item(3, 3).Value = "Intermediate tables for RD3 releases"

# --- entities sheet --------------------------------------------------------
$entities = $wb.Worksheets.Item("entities")
# Row 2 holds the single entity (previously "novelwgs" / "Novel WGS")
$entities.Cells.Item(2, 2).Value = "rd3_portal_release_freeze3"
$entities.Cells.Item(2, 3).Value = "Freeze 3"
$entities.Cells.Item(2, 4).Value = " Staging table for Freeze 3 (2022-03-09)"

# --- attributes sheet -------------------------------------------------------
$attributes = $wb.Worksheets.Item("attributes")
# Column A (entity) for every data row (2-25) referenced the old entity name;
# repoint them all at the renamed entity.
$attributes.Range("A2:A25").Value = "rd3_portal_release_freeze3"
